$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PoolA")

# Sort the data range (A2:H42, with header in row 1) by Date ascending, then by Pool ascending
$sortRange = $ws.Range("A2:H42")
$dateKey = $ws.Range("D2:D42")
$poolKey = $ws.Range("B2:B42")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($dateKey, 0, 1, 0, 0) | Out-Null
$ws.Sort.SortFields.Add($poolKey, 0, 1, 0, 0) | Out-Null

$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.SortMethod = 1
$ws.Sort.Apply()

# The ID column (A) is a manually maintained sequence independent of the sort
# order; restore it to 1..40 after the row reorder.
for ($i = 0; $i -lt 40; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i + 1
}

# Restore view: frozen pane at row1, top-left at A2, active cell G20
$ws.Activate()
$ws.Range("G20").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
